$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "diputados" column (C) is re-curated from a dimension to a measure.
$ws.Range("C2").Value = "iaest-measure:diputados"
$ws.Range("C3").Value = "medida"
$ws.Range("C4").Value = "xsd:int"

# Its mapping file is no longer referenced, so clear that cell entirely.
$ws.Range("C5").Clear()
